$d = $word.ActiveDocument

# The cover letter body currently has 5 paragraphs:
#   1. "I am very interested in applying ..."   (text)
#   2. (empty)
#   3. "I started programming five years ago ..." (text, multiple runs)
#   4. (empty)
#   5. "Thank you so much for reviewing ..."    (text)
#
# The edit removes all the body text, leaving just a single empty
# paragraph before the section properties.

# 1) Remove the first paragraph entirely (its text and its paragraph
#    mark), so the originally-empty second paragraph becomes the first.
$p1 = $d.Paragraphs.Item(1)
$p2 = $d.Paragraphs.Item(2)
$d.Range($p1.Range.Start, $p2.Range.Start).Delete()

# 2) Remove everything from the (now) second paragraph through the end
#    of the document, leaving only the lone empty paragraph behind.
$p2 = $d.Paragraphs.Item(2)
$last = $d.Paragraphs.Item($d.Paragraphs.Count)
$d.Range($p2.Range.Start, $last.Range.End).Delete()
